$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "actual_quantity" / "actual_amount" headers to
# "actual_budget_quantity" / "actual_budget_amount" to match the database.
$ws.Range("M1").Value = "actual_budget_quantity"
$ws.Range("N1").Value = "actual_budget_amount"

# Move the active selection to N1 (matches the saved workbook state).
$ws.Range("N1").Select()
